$wb = $excel.ActiveWorkbook

# "OFF" sheet - row 2 (Home) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 204
$wsOff.Range("C2").Value = 145
$wsOff.Range("D2").Value = 52
$wsOff.Range("E2").Value = 20

# "DEF" sheet - row 2 (Home) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 208
$wsDef.Range("C2").Value = 138
$wsDef.Range("D2").Value = 53
$wsDef.Range("E2").Value = 20
